$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.315.68"
$ws.Range("E2").Value = "  -4.41%  "
$ws.Range("D3").Value = "2.368.19"
$ws.Range("E3").Value = "  -5.86%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'510.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.37%  "
$ws.Range("D6").Value = "'128.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.94%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").Value = "2.375.21"
$ws.Range("E9").Value = "  -5.50%  "
$ws.Range("D10").Value = "'0.0953"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("E12").Value = "  -8.98%  "
$ws.Range("E13").Value = "  -4.97%  "
$ws.Range("D14").Value = "2.789.14"
$ws.Range("E14").Value = "  -5.66%  "
$ws.Range("D15").Value = "56.277.96"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").Value = "'21.44"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("E17").Value = "  -3.85%  "
$ws.Range("D18").Value = "2.360.06"
$ws.Range("E18").Value = "  -5.70%  "
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("E20").Value = "  -4.45%  "
$ws.Range("D21").Value = "'310.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").Value = "'6.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'64.91"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").Value = "'0.390"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.92%  "
$ws.Range("D27").Value = "2.472.48"
$ws.Range("E27").Value = "  -5.70%  "
$ws.Range("E28").Value = "  -5.15%  "
$ws.Range("E29").Value = "  -3.63%  "
$ws.Range("D30").Value = "'174.40"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").Value = "  -6.22%  "
$ws.Range("E34").Value = "  -6.95%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'0.997"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").Value = "'17.62"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("E38").Value = "  -3.92%  "
$ws.Range("E39").Value = "  -6.13%  "
$ws.Range("D40").Value = "'35.60"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").Value = "'0.789"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("D44").Value = "'126.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.64%  "
$ws.Range("E45").Value = "  -4.75%  "
$ws.Range("D46").Value = "'254.68"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.76%  "
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("D48").Value = "'0.0901"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("E49").Value = "  -4.76%  "
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").Value = "'16.58"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.20%  "
